{"js": "// Office.js (Word JavaScript API) implementation of the documented edit.\n//\n// 1) Collapse the \"css\" / \"divTextClass\" / \"buttonClass\" spell-check-flagged\n//    runs into plain text, and flip the \"tableClass\" proofing mark from a\n//    spell-check flag (spellStart/spellEnd) to a grammar flag\n//    (gramStart/gramEnd).\n// 2) Remove the trailing \"Update:\" / \"The \" paragraphs and instead leave a\n//    single trailing space appended (as a new run) to the preceding\n//    \"Design Shell Demo Output.png ...\" paragraph.\n\nconst body = context.document.body;\n\n// --- Part 1: rewrite the \"I have created a styles.css ...\" paragraph -----\nconst oldSentence =\n  \"I have created a styles.css which contains templates for 3 class \" +\n  \"selectors with 3-5 css properties to be filled in and used in HTML. \" +\n  \"The class selectors are to be used for the div class, button class, \" +\n  \"and table class. Its class names are divTextClass, buttonClass, and \" +\n  \"tableClass respectively.\";\n\nconst sentenceHits = body.search(oldSentence, { matchCase: true });\nawait context.sync();\n\nif (sentenceHits.items.length > 0) {\n  const targetRange = sentenceHits.items[0];\n\n  // FlatOPC wrapper so insertOoxml can splice plain runs + proofErr marks\n  // into place without disturbing the paragraph's own identity\n  // (w14:paraId / rsid attributes, etc.) since we are only replacing the\n  // paragraph's inner content, not the <w:p> itself.\n  const replacementOoxml =\n    '<?xml version=\"1.0\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    \"<w:body><w:p>\" +\n    '<w:r><w:t xml:space=\"preserve\">I have created a styles.css which contains templates for 3 class selectors with 3-5 css properties to be filled in and used in HTML. The class selectors are to be used for the div class, button class, and table class. Its class names are divTextClass, buttonClass, and </w:t></w:r>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    \"<w:r><w:t>tableClass</w:t></w:r>\" +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> respectively.</w:t></w:r>' +\n    \"</w:p></w:body></w:document>\" +\n    \"</pkg:xmlData></pkg:part></pkg:package>\";\n\n  targetRange.insertOoxml(replacementOoxml, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- Part 2: drop the \"Update:\" / \"The \" paragraphs, add trailing space ---\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nconst updateIdx = items.findIndex((p) => p.text === \"Update:\");\nconst theIdx = items.findIndex((p) => p.text === \"The \");\nconst demoOutputIdx = items.findIndex(\n  (p) => p.text === \"Design Shell Demo Output.png is a screenshot of the HTML demo output\"\n);\n\nif (demoOutputIdx !== -1) {\n  // Append a standalone run containing a single space to the end of the\n  // \"Design Shell Demo Output.png ...\" paragraph (matches the diff's new\n  // trailing <w:r><w:t xml:space=\"preserve\"> </w:t></w:r>).\n  const demoOutputRange = items[demoOutputIdx].getRange(Word.RangeLocation.end);\n  const spaceRunOoxml =\n    '<?xml version=\"1.0\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    \"<w:body><w:p>\" +\n    '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n    \"</w:p></w:body></w:document>\" +\n    \"</pkg:xmlData></pkg:part></pkg:package>\";\n  demoOutputRange.insertOoxml(spaceRunOoxml, Word.InsertLocation.end);\n  await context.sync();\n}\n\n// Delete \"The \" first (it is after \"Update:\"), then \"Update:\" itself, so\n// indices for the earlier paragraph stay valid while deleting.\nif (theIdx !== -1) {\n  items[theIdx].delete();\n}\nif (updateIdx !== -1) {\n  items[updateIdx].delete();\n}\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) implementation of the documented edit.\n#\n# 1) Collapse the \"css\" / \"divTextClass\" / \"buttonClass\" spell-check-flagged\n#    runs into plain text, and flip the \"tableClass\" proofing mark from a\n#    spell-check flag (spellStart/spellEnd) to a grammar flag\n#    (gramStart/gramEnd).\n# 2) Remove the trailing \"Update:\" / \"The \" paragraphs and instead leave a\n#    single trailing space appended (as a new run) to the preceding\n#    \"Design Shell Demo Output.png ...\" paragraph.\n\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# Part 1: rewrite the \"I have created a styles.css ...\" paragraph\n# ---------------------------------------------------------------------\n$targetPara = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"I have created a styles.css*\") {\n        $targetPara = $p\n        break\n    }\n}\n\nif ($targetPara -ne $null) {\n    # Preserve the paragraph's own identity (w14:paraId) across the\n    # delete+reinsert below - it's the only per-paragraph identifier the\n    # Word object model exposes.\n    $paraId = $targetPara.ParaId\n\n    $r = $targetPara.Range\n    [void]$r.MoveEnd(1, -1)   # exclude the paragraph mark from the range\n    [void]$r.Delete()\n\n    $newParaXml =\n        '<?xml version=\"1.0\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">' +\n        '<w:body><w:p w14:paraId=\"' + $paraId + '\">' +\n        '<w:r><w:t xml:space=\"preserve\">I have created a styles.css which contains templates for 3 class selectors with 3-5 css properties to be filled in and used in HTML. The class selectors are to be used for the div class, button class, and table class. Its class names are divTextClass, buttonClass, and </w:t></w:r>' +\n        '<w:proofErr w:type=\"gramStart\"/>' +\n        '<w:r><w:t>tableClass</w:t></w:r>' +\n        '<w:proofErr w:type=\"gramEnd\"/>' +\n        '<w:r><w:t xml:space=\"preserve\"> respectively.</w:t></w:r>' +\n        '</w:p></w:body></w:document>' +\n        '</pkg:xmlData></pkg:part></pkg:package>'\n\n    [void]$r.InsertXML($newParaXml)\n}\n\n# ---------------------------------------------------------------------\n# Part 2: trailing space on \"Design Shell Demo Output...\" + remove the\n#         \"Update:\" / \"The \" paragraphs that followed it.\n# ---------------------------------------------------------------------\n$demoPara = $null\n$updatePara = $null\n$thePara = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text\n    if ($t -eq \"Design Shell Demo Output.png is a screenshot of the HTML demo output`r\") {\n        $demoPara = $p\n    } elseif ($t -eq \"Update:`r\") {\n        $updatePara = $p\n    } elseif ($t -eq \"The `r\") {\n        $thePara = $p\n    }\n}\n\nif ($demoPara -ne $null) {\n    $r = $demoPara.Range\n    [void]$r.MoveEnd(1, -1)   # exclude the paragraph mark, keep the existing run(s) intact\n    $spaceRunXml =\n        '<?xml version=\"1.0\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body><w:p><w:r><w:t xml:space=\"preserve\"> </w:t></w:r></w:p></w:body></w:document>' +\n        '</pkg:xmlData></pkg:part></pkg:package>'\n    [void]$r.InsertXML($spaceRunXml)\n}\n\n# Delete \"The \" before \"Update:\" so paragraph references above stay valid.\nif ($thePara -ne $null) {\n    [void]$thePara.Range.Delete()\n}\nif ($updatePara -ne $null) {\n    [void]$updatePara.Range.Delete()\n}\n"}
